$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix Primary Objective text (typo / missing apostrophe) ---
$ws.Range("E2").Value = "Evaluate the efficacy and safety of transdermal xanomeline, 50cm2 and 75cm2, and placebo in subjects with mild to moderate Alzheimer's disease."

# --- Fix Primary Outcome Measure casing ---
$ws.Range("E3").Value = "ADAS-Cog"

# --- Align visit code column (G) with visit label column (D); drop stray lowercase duplicates ---
# Row 18 keeps its existing highlight style, only the text changes.
$ws.Range("G18").Value = "Wk16"

# Rows 19-25 lose their highlight style entirely (Clear resets formatting+contents), then get the corrected text.
$ws.Range("G19").Clear()
$ws.Range("G19").Value = "Wk20"

$ws.Range("G20").Clear()
$ws.Range("G20").Value = "Wk24"

$ws.Range("G21").Clear()
$ws.Range("G21").Value = "Wk26"

$ws.Range("G22").Clear()
$ws.Range("G22").Value = "Wk2"

$ws.Range("G23").Clear()
$ws.Range("G23").Value = "Wk4"

$ws.Range("G24").Clear()
$ws.Range("G24").Value = "Wk6"

$ws.Range("G25").Clear()
$ws.Range("G25").Value = "Wk8"

# --- Start conversion of AE: add new studytype/phase parameter rows ---
$ws.Cells.Item(26, 1).Value = "CDISCPILOT01"
$ws.Cells.Item(26, 2).Value = 1
$ws.Cells.Item(26, 3).Value = "studytype"
$ws.Cells.Item(26, 4).Value = "Study Type"
$ws.Cells.Item(26, 5).Value = "INT"

$ws.Cells.Item(27, 1).Value = "CDISCPILOT01"
$ws.Cells.Item(27, 2).Value = 1
$ws.Cells.Item(27, 3).Value = "phase"
$ws.Cells.Item(27, 4).Value = "Trial Phase"
$ws.Cells.Item(27, 5).Value = "PHASE2"

# --- Restore the active selection to match the saved workbook state ---
$ws.Range("E13").Select()
